# Update numeric "想去人数" (F column) values on the "展览" and "全部类型"
# worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 610
    $ws.Range("F23").Value = 476
    $ws.Range("F24").Value = 5045
    $ws.Range("F25").Value = 247
}
